$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated run_time, max_er, and iteration columns (iter 0..19) for rows 2-11
# per data_set "gr50_02" — low-input values for other-language run (see commit message).
$updates = @{
    "C2" = 0.8209939002990723
    "E2" = 513.1569947413573
    "F2" = 0.01699441986346073
    "G2" = 0.01429048309485674
    "H2" = 0.01358705561238463
    "I2" = 0.01290113526781863
    "J2" = 0.01253727704768912
    "K2" = 0.01184749830946254
    "L2" = 0.01182968960312058
    "M2" = 0.01166861166523757
    "N2" = 0.01097272640344791
    "O2" = 0.01084067661433536
    "P2" = 0.01062463163492577
    "Q2" = 0.01047415713280885
    "R2" = 0.01044886964900599
    "S2" = 0.01027126715240739
    "T2" = 0.01024634083605576
    "U2" = 0.01011890419877931
    "V2" = 0.01010612286048982
    "W2" = 0.01005232967494557
    "X2" = 0.01001116206992961
    "Y2" = 0.01000306032634224
    "C3" = 0.8009719848632812
    "E3" = 529.2264602940522
    "F3" = 0.01728406644392292
    "G3" = 0.01462761415928338
    "H3" = 0.01361299423191866
    "I3" = 0.01317889976324923
    "J3" = 0.01273481341914853
    "K3" = 0.01218946135084497
    "L3" = 0.01178379870465622
    "M3" = 0.011145747795459
    "N3" = 0.011145747795459
    "O3" = 0.01109363968650316
    "P3" = 0.01081791440071315
    "Q3" = 0.01070632456863787
    "R3" = 0.01058296904372875
    "S3" = 0.01057170513358791
    "T3" = 0.01056611995870828
    "U3" = 0.01056594678741915
    "V3" = 0.01048808092181851
    "W3" = 0.01036042678632979
    "X3" = 0.01031910837234223
    "Y3" = 0.0103163052688899
    "C4" = 0.7680397033691406
    "E4" = 535.3558655305696
    "F4" = 0.01725169053182676
    "G4" = 0.01463952231815306
    "H4" = 0.01351134402526451
    "I4" = 0.01219645953439066
    "J4" = 0.01214171104518077
    "K4" = 0.01175875890033145
    "L4" = 0.01138051294299685
    "M4" = 0.01105559669381131
    "N4" = 0.01105559669381131
    "O4" = 0.01086157715261722
    "P4" = 0.01086157715261722
    "Q4" = 0.01086157715261722
    "R4" = 0.01086157715261722
    "S4" = 0.01078149647626602
    "T4" = 0.01078149647626602
    "U4" = 0.01068848261947878
    "V4" = 0.01061436443234329
    "W4" = 0.01051390503829733
    "X4" = 0.01045329614920757
    "Y4" = 0.01043578685244775
    "C5" = 0.8279554843902588
    "E5" = 506.8226562451946
    "F5" = 0.01739592327102118
    "G5" = 0.01437089899969145
    "H5" = 0.01297226798244558
    "I5" = 0.01232688841922479
    "J5" = 0.01158717129006685
    "K5" = 0.01132395312355612
    "L5" = 0.01115751860601796
    "M5" = 0.01080705533529358
    "N5" = 0.01048582095318901
    "O5" = 0.01048582095318901
    "P5" = 0.01033673981795818
    "Q5" = 0.0101884951688873
    "R5" = 0.01011572108657935
    "S5" = 0.01008878121370462
    "T5" = 0.01008878121370462
    "U5" = 0.009979098135520658
    "V5" = 0.00991270249561599
    "W5" = 0.00991270249561599
    "X5" = 0.009882697914407023
    "Y5" = 0.009879583942401451
    "C6" = 0.875
    "E6" = 513.9236668243957
    "F6" = 0.01658891819224734
    "G6" = 0.01428184628700338
    "H6" = 0.01357345332731523
    "I6" = 0.01220641253182412
    "J6" = 0.01219209698615055
    "K6" = 0.01192549183119938
    "L6" = 0.01174669047412513
    "M6" = 0.01124755803578691
    "N6" = 0.01073798235681318
    "O6" = 0.01062757546358012
    "P6" = 0.01042264564235825
    "Q6" = 0.01042264564235825
    "R6" = 0.0103055550362119
    "S6" = 0.01024773311152531
    "T6" = 0.01013419737808607
    "U6" = 0.0100878230143616
    "V6" = 0.0100878230143616
    "W6" = 0.01005084321099509
    "X6" = 0.01004064426738787
    "Y6" = 0.01001800520125527
    "C7" = 0.7600007057189941
    "E7" = 528.5294521073974
    "F7" = 0.0172961155805952
    "G7" = 0.01491511490831735
    "H7" = 0.01372559013243485
    "I7" = 0.01297079353106389
    "J7" = 0.01241895861988799
    "K7" = 0.01210908021703531
    "L7" = 0.01152102895026909
    "M7" = 0.0109656523657843
    "N7" = 0.0109656523657843
    "O7" = 0.01069360795244022
    "P7" = 0.01069360795244022
    "Q7" = 0.01069360795244022
    "R7" = 0.01068739181830065
    "S7" = 0.0105725677969979
    "T7" = 0.01054886029687172
    "U7" = 0.01050800329762675
    "V7" = 0.01044515654690903
    "W7" = 0.01044515654690903
    "X7" = 0.01031086118057032
    "Y7" = 0.01030271836466662
    "C8" = 0.8870003223419189
    "E8" = 521.7736329281743
    "F8" = 0.01672757691776776
    "G8" = 0.0144581959140737
    "H8" = 0.01328793770987029
    "I8" = 0.01227406682538534
    "J8" = 0.01203587357654638
    "K8" = 0.01164943005893448
    "L8" = 0.01105700934339192
    "M8" = 0.01093242022550481
    "N8" = 0.01065785660341695
    "O8" = 0.01065785660341695
    "P8" = 0.01063151319753667
    "Q8" = 0.01060661961941054
    "R8" = 0.01040639895621504
    "S8" = 0.01036598914802359
    "T8" = 0.01034217171349137
    "U8" = 0.01024984276492285
    "V8" = 0.01024984276492285
    "W8" = 0.01020094281197773
    "X8" = 0.01017921737004449
    "Y8" = 0.01017102598300534
    "C9" = 0.8649981021881104
    "E9" = 507.0114692100906
    "F9" = 0.01653044269179504
    "G9" = 0.01407110700368684
    "H9" = 0.01300150399856541
    "I9" = 0.01265988461497353
    "J9" = 0.01178383600954401
    "K9" = 0.01121579386339636
    "L9" = 0.01093677142707396
    "M9" = 0.01093677142707396
    "N9" = 0.01077782402845285
    "O9" = 0.01057963871002609
    "P9" = 0.01023555351962147
    "Q9" = 0.01023555351962147
    "R9" = 0.01023221486802197
    "S9" = 0.01023221486802197
    "T9" = 0.01011757219620968
    "U9" = 0.01002423794185428
    "V9" = 0.01000278172784588
    "W9" = 0.009907610390024205
    "X9" = 0.009890286066042354
    "Y9" = 0.009883264507019308
    "C10" = 0.7659986019134521
    "E10" = 511.3681778331811
    "F10" = 0.0173752586922241
    "G10" = 0.01399414797609602
    "H10" = 0.01317817488472904
    "I10" = 0.01252350209954949
    "J10" = 0.01217424239386712
    "K10" = 0.01162496737040173
    "L10" = 0.0111660571181352
    "M10" = 0.0111660571181352
    "N10" = 0.01074514116767926
    "O10" = 0.01056091676816136
    "P10" = 0.01056091676816136
    "Q10" = 0.01037861553111256
    "R10" = 0.01034918838048368
    "S10" = 0.01018144039217588
    "T10" = 0.01018144039217588
    "U10" = 0.01016620890045812
    "V10" = 0.01002236082127209
    "W10" = 0.01002236082127209
    "X10" = 0.009997035171590889
    "Y10" = 0.009968190601036667
    "C11" = 0.7499940395355225
    "E11" = 512.1807076651858
    "F11" = 0.01698559749220471
    "G11" = 0.01485421660268369
    "H11" = 0.01332602461891194
    "I11" = 0.01255697708627753
    "J11" = 0.01138068118170998
    "K11" = 0.01132450853997211
    "L11" = 0.01102724018559378
    "M11" = 0.01075052625084127
    "N11" = 0.01063384181052271
    "O11" = 0.01049569421097855
    "P11" = 0.01049569421097855
    "Q11" = 0.01046473542508529
    "R11" = 0.01039481343843228
    "S11" = 0.0103480411408149
    "T11" = 0.01021055383813355
    "U11" = 0.01021055383813355
    "V11" = 0.01007462934904373
    "W11" = 0.01006729871989669
    "X11" = 0.01001202640107415
    "Y11" = 0.009984029389184905
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
